# Refresh the cryptocurrency price/volume snapshot in columns D (Price) and E
# (Volume 1h change). Values that would otherwise be auto-parsed as numbers by
# Excel (losing their original text formatting, e.g. trailing zeros like
# "23.50") are written with a leading apostrophe so they stay plain text,
# exactly like the other price strings already on the sheet (e.g. "1.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.967.34"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.647.29"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D5").Value = "'213.75"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'0.527"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'23.50"
$ws.Range("E8").Value = "  +2.87%  "
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "'0.0615"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "'0.0872"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "1.881.13"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Value = "1.643.14"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "'0.566"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "'65.65"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "27.999.65"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "'232.07"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  +5.10%  "
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").Value = "'152.27"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "'6.92"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'15.78"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'1.20"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("D33").Value = "1.444.26"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'0.890"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "'0.935"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("D40").Value = "'0.559"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'69.25"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +5.56%  "
$ws.Range("D46").Value = "'5.42"
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "1.789.56"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "'89.16"
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("E51").Value = "  +0.09%  "
